$d = $word.ActiveDocument

# Paragraph 2 already reads: " Feel so boring1" + "I am bob, every day I need to write a report, "
# We need to split paragraph 1 ("I am bob, every day I need to write a report, ")
# into two paragraphs by inserting a duplicate of that pattern right before paragraph 2,
# i.e. insert a new paragraph, right after paragraph 1, that reads:
#   " Feel so boring1" + "I am bob, every day I need to write a report, "
$srcPara = $d.Paragraphs.Item(2)
$srcFormatted = $srcPara.Range.FormattedText

$insertPoint = $d.Range($srcPara.Range.Start, $srcPara.Range.Start)
$insertPoint.FormattedText = $srcFormatted
